# Bump the fixed "Date Placeholder" text stamped on the slide master and
# every slide layout from 3/21/21 to 3/24/21.

$p = $ppt.ActivePresentation

$oldDate = "3/21/21"
$newDate = "3/24/21"

function Update-DateShapes {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText -and $tf.TextRange.Text -eq $oldDate) {
                $tf.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout belonging to the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DateShapes $layouts.Item($l).Shapes
}

# Cover any additional masters/designs in the deck as well.
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $master = $p.Designs.Item($d).SlideMaster
    Update-DateShapes $master.Shapes
    $dLayouts = $master.CustomLayouts
    for ($l = 1; $l -le $dLayouts.Count; $l++) {
        Update-DateShapes $dLayouts.Item($l).Shapes
    }
}
